$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the Area / Atotal columns, plus the small
# Atotal/Qtotal summary pair in J1:K1.
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Re-key the "segment" (D) formula across D3:D9 into its own shared-formula
# group (values are unchanged, this mirrors the source file's layout).
$ws.Range("D3:D9").Formula = "=(A3/100+(A4/100-A3/100)/2)"

# New Area (G) column.
# Row 2 is measured from 0 (no previous depth reading).
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
# Row 3 follows the normal incremental-area pattern.
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
# Rows 4-15 are one continuous fill (shared formula) of the same pattern.
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Atotal (H2) sums the new Area column.
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Summary cross references next to the table.
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

$ws.Range("J2:K2").Select() | Out-Null
